# Apply crypto price/volume updates (commit: "Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates -------------------------------------------
# These values look numeric (e.g. "0.4350", "15.78") so a plain .Value =
# assignment would let Excel re-parse them as doubles and silently drop
# trailing zeros / switch to scientific notation. Force text entry via the
# "@" number format, then ClearFormats() so the cell keeps the workbook's
# original (unstyled) formatting -- only the stored text should change.
$priceUpdates = [ordered]@{
    'D2' = '28.004.57'
    'D3' = '1.859.43'
    'D5' = '317.17'
    'D7' = '0.4350'
    'D8' = '0.3674'
    'D9' = '0.07483'
    'D10' = '0.9350'
    'D11' = '21.27'
    'D12' = '1.860.78'
    'D13' = '6.688'
    'D14' = '5.416'
    'D15' = '0.06898'
    'D16' = '1.006'
    'D18' = '0.000008995'
    'D20' = '15.78'
    'D21' = '28.015.98'
    'D22' = '5.106'
    'D23' = '10.81'
    'D24' = '2.111.52'
    'D25' = '2.019'
    'D26' = '154.08'
    'D27' = '18.34'
    'D30' = '1.726'
    'D31' = '0.08961'
    'D32' = '0.7966'
    'D34' = '3.039'
    'D35' = '1.169'
    'D36' = '1.003'
    'D38' = '0.05411'
    'D39' = '0.01949'
    'D40' = '2.926'
    'D41' = '0.5220'
    'D42' = '6.979'
    'D43' = '0.1676'
    'D44' = '8.692'
    'D45' = '0.06708'
    'D46' = '0.4856'
    'D47' = '10.58'
    'D48' = '106.67'
    'D50' = '1.901'
    'D51' = '1.665'
}
foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.ClearFormats()
}

# --- Column E ("Volume(1h)") updates ---------------------------------------
# Padded with literal leading/trailing spaces, so Excel never mistakes these
# for numeric/percentage values -- a direct .Value assignment is safe.
$volumeUpdates = [ordered]@{
    'E2' = '  -3.30%  '
    'E3' = '  -2.72%  '
    'E4' = '  +0.13%  '
    'E5' = '  -2.37%  '
    'E6' = '  +0.16%  '
    'E7' = '  -5.19%  '
    'E9' = '  -3.10%  '
    'E10' = '  -4.70%  '
    'E11' = '  -3.70%  '
    'E12' = '  -2.11%  '
    'E13' = '  -3.67%  '
    'E14' = '  -4.37%  '
    'E15' = '  -2.10%  '
    'E16' = '  +0.13%  '
    'E17' = '  -3.03%  '
    'E18' = '  -4.91%  '
    'E19' = '  +0.26%  '
    'E20' = '  -5.54%  '
    'E22' = '  -4.01%  '
    'E23' = '  -0.78%  '
    'E24' = '  -0.98%  '
    'E25' = '  -3.57%  '
    'E26' = '  -2.76%  '
    'E27' = '  -3.79%  '
    'E28' = '  -5.47%  '
    'E29' = '  -3.80%  '
    'E30' = '  -6.67%  '
    'E31' = '  -3.58%  '
    'E32' = '  -8.28%  '
    'E33' = '  -5.32%  '
    'E34' = '  +0.33%  '
    'E35' = '  -6.38%  '
    'E36' = '  +0.13%  '
    'E37' = '  -3.24%  '
    'E38' = '  -5.50%  '
    'E39' = '  -4.57%  '
    'E40' = '  +3.20%  '
    'E41' = '  -5.08%  '
    'E42' = '  -5.64%  '
    'E43' = '  -4.37%  '
    'E44' = '  -6.71%  '
    'E45' = '  -2.32%  '
    'E46' = '  -6.37%  '
    'E47' = '  -5.89%  '
    'E48' = '  -3.53%  '
    'E49' = '  +0.08%  '
    'E50' = '  -7.64%  '
    'E51' = '  -6.50%  '
}
foreach ($ref in $volumeUpdates.Keys) {
    $ws.Range($ref).Value = $volumeUpdates[$ref]
}
